# Update the NATMI ligand-receptor pair table (Cx3cl1-Itgav) with newly
# computed TPM-derived statistics.
#
# The worksheet has one row per (Sending cluster, Target cluster) pair
# (4 sending clusters x 4 target clusters = 16 data rows, rows 2..17).
# Columns:
#   A Sending cluster          B Ligand symbol          C Receptor symbol
#   D Target cluster
#   E Ligand-expressing cells  F Ligand detection rate
#   G Ligand average expr.     H Ligand total expr.
#   I Ligand avg specificity   J Ligand total specificity
#   K Receptor-expressing cells  L Receptor detection rate
#   M Receptor average expr.   N Receptor total expr.
#   O Receptor avg specificity P Receptor total specificity
#   Q Edge average expr weight R Edge total expr weight
#   S Edge average expr specificity T Edge total expr specificity
#
# The ligand-side statistics (E:J) only depend on the Sending cluster (A)
# and the receptor-side statistics (K:P) only depend on the Target cluster
# (D). The edge columns (Q:T) are simply the products of the matching
# ligand/receptor average & total values and their specificities:
#   Q = G * M   R = H * N   S = I * O   T = J * P

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ligand-side values (E, F, G, H, I, J) keyed by Sending cluster.
$ligand = @{
    "ECs"           = @(3, 1,                 11.64940266666667, 34.948208,          0.4844115508208772,  0.4844115508208771)
    "FAPs"          = @(3, 1,                 10.423773,         31.271319,          0.4334467773856777,  0.4334467773856777)
    "MuSCs"         = @(3, 1,                 1.677245,          5.031734999999999,  0.06974407828492055, 0.06974407828492052)
    "Resolving-Mac" = @(2, 0.6666666666666666,0.2981443333333333,0.894433,           0.01239759350852466, 0.01239759350852466)
}

# New receptor-side values (K, L, M, N, O, P) keyed by Target cluster.
$receptor = @{
    "ECs"           = @(3, 1, 8.820647333333334, 26.461942,   0.06415146660411865, 0.06415146660411865)
    "FAPs"          = @(3, 1, 54.711535,          164.134605,  0.3979101621202897,  0.3979101621202898)
    "MuSCs"         = @(3, 1, 21.90816333333333,  65.72449,    0.1593353362087987,  0.1593353362087987)
    "Resolving-Mac" = @(3, 1, 52.056859,          156.170577,  0.3786030350667928,  0.3786030350667929)
}

$lastRow = 17
for ($r = 2; $r -le $lastRow; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value()
    $target  = $ws.Cells.Item($r, 4).Value()

    $l = $ligand[$sending]
    $rc = $receptor[$target]

    # E..J : ligand stats
    $ws.Cells.Item($r, 5).Value  = $l[0]
    $ws.Cells.Item($r, 6).Value  = $l[1]
    $ws.Cells.Item($r, 7).Value  = $l[2]
    $ws.Cells.Item($r, 8).Value  = $l[3]
    $ws.Cells.Item($r, 9).Value  = $l[4]
    $ws.Cells.Item($r, 10).Value = $l[5]

    # K..P : receptor stats
    $ws.Cells.Item($r, 11).Value = $rc[0]
    $ws.Cells.Item($r, 12).Value = $rc[1]
    $ws.Cells.Item($r, 13).Value = $rc[2]
    $ws.Cells.Item($r, 14).Value = $rc[3]
    $ws.Cells.Item($r, 15).Value = $rc[4]
    $ws.Cells.Item($r, 16).Value = $rc[5]

    # Q..T : edge stats = product of ligand & receptor average/total/specificity
    $ws.Cells.Item($r, 17).Value = $l[2] * $rc[2]
    $ws.Cells.Item($r, 18).Value = $l[3] * $rc[3]
    $ws.Cells.Item($r, 19).Value = $l[4] * $rc[4]
    $ws.Cells.Item($r, 20).Value = $l[5] * $rc[5]
}
